# Generate Report for handback
# Updates the "9fc327a9-32f4-4496-87fb-22743577c874" row from
# "Ready for handoff" to "Handed back: in sync with en-US" across all
# sheets, and records new "Latest Handback DateTime" values for the
# zh-cn and de-de locale sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: Row 3 is the 9fc327a9 file, columns B (zh-cn) and C (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: Row 3 is the 9fc327a9 file
#     B = Status, G = Latest Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-01-18 12:25:26"

# --- de-de sheet: Row 3 is the 9fc327a9 file
#     B = Status, G = Latest Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-01-18 12:25:43"
